$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 37,13

$arr[0,0] = 535
$arr[0,1] = 'Sunday, Jan 15'
$arr[0,2] = '7:52 AM'
$arr[0,3] = 'TO1930'
$arr[0,4] = 'Paris'
$arr[0,5] = '(ORY)'
$arr[0,6] = 'Transavia '
$arr[0,7] = 'B738'
$arr[0,8] = '(F-HTVX)'
$arr[0,9] = '8:21 AM'
$arr[0,11] = '0 hours, 29 minutes'
$arr[1,0] = 536
$arr[1,1] = 'Sunday, Jan 15'
$arr[1,2] = '8:10 AM'
$arr[1,3] = 'FR4425'
$arr[1,4] = 'Cagliari'
$arr[1,5] = '(CAG)'
$arr[1,6] = 'Ryanair '
$arr[1,7] = 'B738'
$arr[1,8] = '(9H-QCV)'
$arr[1,9] = '7:59 AM'
$arr[1,11] = '0 hours, -11 minutes'
$arr[2,0] = 537
$arr[2,1] = 'Sunday, Jan 15'
$arr[2,2] = '8:10 AM'
$arr[2,3] = 'LO3903'
$arr[2,4] = 'Warsaw'
$arr[2,5] = '(WAW)'
$arr[2,6] = 'LOT (Star Alliance Livery) '
$arr[2,7] = 'E75S'
$arr[2,8] = '(SP-LIO)'
$arr[2,9] = '8:05 AM'
$arr[2,11] = '0 hours, -5 minutes'
$arr[3,0] = 538
$arr[3,1] = 'Sunday, Jan 15'
$arr[3,2] = '8:45 AM'
$arr[3,3] = 'LX1370'
$arr[3,4] = 'Zurich'
$arr[3,5] = '(ZRH)'
$arr[3,6] = 'Swiss '
$arr[3,7] = 'BCS3'
$arr[3,8] = '(HB-JCF)'
$arr[3,9] = '8:34 AM'
$arr[3,11] = '0 hours, -11 minutes'
$arr[4,0] = 539
$arr[4,1] = 'Sunday, Jan 15'
$arr[4,2] = '9:20 AM'
$arr[4,3] = 'DY1038'
$arr[4,4] = 'Trondheim'
$arr[4,5] = '(TRD)'
$arr[4,6] = 'Norwegian '
$arr[4,7] = 'B738'
$arr[4,8] = '(LN-ENT)'
$arr[4,9] = '9:15 AM'
$arr[4,11] = '0 hours, -5 minutes'
$arr[5,0] = 540
$arr[5,1] = 'Sunday, Jan 15'
$arr[5,2] = '9:30 AM'
$arr[5,3] = 'FR2362'
$arr[5,4] = 'London'
$arr[5,5] = '(STN)'
$arr[5,6] = 'Ryanair '
$arr[5,7] = 'B738'
$arr[5,8] = '(EI-DWA)'
$arr[5,9] = '9:10 AM'
$arr[5,11] = '0 hours, -20 minutes'
$arr[6,0] = 541
$arr[6,1] = 'Sunday, Jan 15'
$arr[6,2] = '9:40 AM'
$arr[6,3] = 'DY1028'
$arr[6,4] = 'Bergen'
$arr[6,5] = '(BGO)'
$arr[6,6] = 'Norwegian '
$arr[6,7] = 'B738'
$arr[6,8] = '(SE-RPD)'
$arr[6,9] = '9:25 AM'
$arr[6,11] = '0 hours, -15 minutes'
$arr[7,0] = 542
$arr[7,1] = 'Sunday, Jan 15'
$arr[7,2] = '10:20 AM'
$arr[7,3] = 'D84901'
$arr[7,4] = 'Stockholm'
$arr[7,5] = '(ARN)'
$arr[7,6] = 'Norwegian '
$arr[7,7] = 'B738'
$arr[7,8] = '(LN-NGD)'
$arr[7,9] = '10:08 AM'
$arr[7,11] = '0 hours, -12 minutes'
$arr[8,0] = 543
$arr[8,1] = 'Sunday, Jan 15'
$arr[8,2] = '10:30 AM'
$arr[8,3] = 'FR1312'
$arr[8,4] = 'Seville'
$arr[8,5] = '(SVQ)'
$arr[8,6] = 'Ryanair '
$arr[8,7] = 'B738'
$arr[8,8] = '(EI-ENN)'
$arr[8,9] = '10:10 AM'
$arr[8,11] = '0 hours, -20 minutes'
$arr[9,0] = 544
$arr[9,1] = 'Sunday, Jan 15'
$arr[9,2] = '10:30 AM'
$arr[9,3] = 'FR6249'
$arr[9,4] = 'Manchester'
$arr[9,5] = '(MAN)'
$arr[9,6] = 'Ryanair '
$arr[9,7] = 'B38M'
$arr[9,8] = '(EI-HGE)'
$arr[9,9] = '10:17 AM'
$arr[9,11] = '0 hours, -13 minutes'
$arr[10,0] = 545
$arr[10,1] = 'Sunday, Jan 15'
$arr[10,2] = '10:55 AM'
$arr[10,3] = 'FR6257'
$arr[10,4] = 'Stockholm'
$arr[10,5] = '(ARN)'
$arr[10,6] = 'Ryanair '
$arr[10,7] = 'B38M'
$arr[10,8] = '(SP-RZK)'
$arr[10,9] = '10:47 AM'
$arr[10,11] = '0 hours, -8 minutes'
$arr[11,0] = 546
$arr[11,1] = 'Sunday, Jan 15'
$arr[11,2] = '11:05 AM'
$arr[11,3] = 'DY1040'
$arr[11,4] = 'Oslo'
$arr[11,5] = '(OSL)'
$arr[11,6] = 'Norwegian (Jan Baalsrud Livery) '
$arr[11,7] = 'B738'
$arr[11,8] = '(LN-ENR)'
$arr[11,9] = '11:04 AM'
$arr[11,11] = '0 hours, -1 minutes'
$arr[12,0] = 547
$arr[12,1] = 'Sunday, Jan 15'
$arr[12,2] = '11:05 AM'
$arr[12,3] = 'FR6331'
$arr[12,4] = 'Thessaloniki'
$arr[12,5] = '(SKG)'
$arr[12,6] = 'Ryanair '
$arr[12,7] = 'B738'
$arr[12,8] = '(EI-DCN)'
$arr[12,9] = '10:51 AM'
$arr[12,11] = '0 hours, -14 minutes'
$arr[13,0] = 548
$arr[13,1] = 'Sunday, Jan 15'
$arr[13,2] = '11:10 AM'
$arr[13,3] = 'FZ1787'
$arr[13,4] = 'Dubai'
$arr[13,5] = '(DXB)'
$arr[13,6] = 'flydubai '
$arr[13,7] = 'B38M'
$arr[13,8] = '(A6-FMI)'
$arr[13,9] = '11:14 AM'
$arr[13,11] = '0 hours, 4 minutes'
$arr[14,0] = 549
$arr[14,1] = 'Sunday, Jan 15'
$arr[14,2] = '11:15 AM'
$arr[14,3] = 'W65068'
$arr[14,4] = 'Milan'
$arr[14,5] = '(MXP)'
$arr[14,6] = 'Wizz Air '
$arr[14,7] = 'A321'
$arr[14,8] = '(HA-LXO)'
$arr[14,9] = '10:54 AM'
$arr[14,11] = '0 hours, -21 minutes'
$arr[15,0] = 550
$arr[15,1] = 'Sunday, Jan 15'
$arr[15,2] = '11:25 AM'
$arr[15,3] = 'OS597'
$arr[15,4] = 'Vienna'
$arr[15,5] = '(VIE)'
$arr[15,6] = 'Austrian Airlines '
$arr[15,7] = 'E195'
$arr[15,8] = '(OE-LWM)'
$arr[15,9] = '11:08 AM'
$arr[15,11] = '0 hours, -17 minutes'
$arr[16,0] = 551
$arr[16,1] = 'Sunday, Jan 15'
$arr[16,2] = '11:30 AM'
$arr[16,3] = 'BA872'
$arr[16,4] = 'London'
$arr[16,5] = '(LHR)'
$arr[16,6] = 'British Airways '
$arr[16,7] = 'A320'
$arr[16,8] = '(G-EUUD)'
$arr[16,9] = '11:02 AM'
$arr[16,11] = '0 hours, -28 minutes'
$arr[17,0] = 552
$arr[17,1] = 'Sunday, Jan 15'
$arr[17,2] = '11:30 AM'
$arr[17,3] = 'LO3907'
$arr[17,4] = 'Warsaw'
$arr[17,5] = '(WAW)'
$arr[17,6] = 'LOT '
$arr[17,7] = 'E190'
$arr[17,8] = '(SP-LMD)'
$arr[17,9] = '11:22 AM'
$arr[17,11] = '0 hours, -8 minutes'
$arr[18,0] = 553
$arr[18,1] = 'Sunday, Jan 15'
$arr[18,2] = '11:35 AM'
$arr[18,3] = 'FR3510'
$arr[18,4] = 'Milan'
$arr[18,5] = '(BGY)'
$arr[18,6] = 'Buzz '
$arr[18,7] = 'B38M'
$arr[18,8] = '(SP-RZH)'
$arr[18,9] = '11:32 AM'
$arr[18,11] = '0 hours, -3 minutes'
$arr[19,0] = 554
$arr[19,1] = 'Sunday, Jan 15'
$arr[19,2] = '12:05 PM'
$arr[19,3] = 'W65060'
$arr[19,4] = 'Catania'
$arr[19,5] = '(CTA)'
$arr[19,6] = 'Wizz Air '
$arr[19,7] = 'A21N'
$arr[19,8] = '(HA-LZI)'
$arr[19,9] = '11:42 AM'
$arr[19,11] = '0 hours, -23 minutes'
$arr[20,0] = 555
$arr[20,1] = 'Sunday, Jan 15'
$arr[20,2] = '12:10 PM'
$arr[20,3] = 'FR2727'
$arr[20,4] = 'Catania'
$arr[20,5] = '(CTA)'
$arr[20,6] = 'Ryanair '
$arr[20,7] = 'B738'
$arr[20,8] = '(SP-RKB)'
$arr[20,9] = '11:54 AM'
$arr[20,11] = '0 hours, -16 minutes'
$arr[21,0] = 556
$arr[21,1] = 'Sunday, Jan 15'
$arr[21,2] = '12:45 PM'
$arr[21,3] = 'LH1620'
$arr[21,4] = 'Munich'
$arr[21,5] = '(MUC)'
$arr[21,6] = 'Lufthansa '
$arr[21,7] = 'A320'
$arr[21,8] = '(D-AIWB)'
$arr[21,9] = '12:42 PM'
$arr[21,11] = '0 hours, -3 minutes'
$arr[22,0] = 557
$arr[22,1] = 'Sunday, Jan 15'
$arr[22,2] = '12:55 PM'
$arr[22,3] = 'FR9662'
$arr[22,4] = 'Rome'
$arr[22,5] = '(CIA)'
$arr[22,6] = 'Ryanair '
$arr[22,7] = 'B738'
$arr[22,8] = '(9H-QEA)'
$arr[22,9] = '12:39 PM'
$arr[22,11] = '0 hours, -16 minutes'
$arr[23,0] = 558
$arr[23,1] = 'Sunday, Jan 15'
$arr[23,2] = '1:00 PM'
$arr[23,3] = 'FR6159'
$arr[23,4] = 'Gothenburg'
$arr[23,5] = '(GOT)'
$arr[23,6] = 'Ryanair '
$arr[23,7] = 'B738'
$arr[23,8] = '(EI-DWY)'
$arr[23,9] = '12:49 PM'
$arr[23,11] = '0 hours, -11 minutes'
$arr[24,0] = 559
$arr[24,1] = 'Sunday, Jan 15'
$arr[24,2] = '1:15 PM'
$arr[24,3] = 'FR2023'
$arr[24,4] = 'Dublin'
$arr[24,5] = '(DUB)'
$arr[24,6] = 'Buzz '
$arr[24,7] = 'B38M'
$arr[24,8] = '(SP-RZD)'
$arr[24,9] = '1:07 PM'
$arr[24,11] = '0 hours, -8 minutes'
$arr[25,0] = 560
$arr[25,1] = 'Sunday, Jan 15'
$arr[25,2] = '1:25 PM'
$arr[25,3] = 'FR6247'
$arr[25,4] = 'Pafos'
$arr[25,5] = '(PFO)'
$arr[25,6] = 'Buzz '
$arr[25,7] = 'B38M'
$arr[25,8] = '(SP-RZF)'
$arr[25,9] = '1:33 PM'
$arr[25,11] = '0 hours, 8 minutes'
$arr[26,0] = 561
$arr[26,1] = 'Sunday, Jan 15'
$arr[26,2] = '1:30 PM'
$arr[26,3] = 'U23815'
$arr[26,4] = 'Paris'
$arr[26,5] = '(CDG)'
$arr[26,6] = 'easyJet '
$arr[26,7] = 'A320'
$arr[26,8] = '(OE-INB)'
$arr[26,9] = '1:28 PM'
$arr[26,11] = '0 hours, -2 minutes'
$arr[27,0] = 562
$arr[27,1] = 'Sunday, Jan 15'
$arr[27,2] = '1:45 PM'
$arr[27,3] = 'FR3798'
$arr[27,4] = 'Bologna'
$arr[27,5] = '(BLQ)'
$arr[27,6] = 'Ryanair '
$arr[27,7] = 'B738'
$arr[27,8] = '(9H-QAM)'
$arr[27,9] = '2:08 PM'
$arr[27,11] = '0 hours, 23 minutes'
$arr[28,0] = 563
$arr[28,1] = 'Sunday, Jan 15'
$arr[28,2] = '2:20 PM'
$arr[28,3] = 'LO3905'
$arr[28,4] = 'Warsaw'
$arr[28,5] = '(WAW)'
$arr[28,6] = 'LOT '
$arr[28,7] = 'E195'
$arr[28,8] = '(SP-LNP)'
$arr[28,9] = '2:13 PM'
$arr[28,11] = '0 hours, -7 minutes'
$arr[29,0] = 564
$arr[29,1] = 'Sunday, Jan 15'
$arr[29,2] = '3:05 PM'
$arr[29,3] = 'FR6235'
$arr[29,4] = 'Copenhagen'
$arr[29,5] = '(CPH)'
$arr[29,6] = 'Ryanair '
$arr[29,7] = 'B38M'
$arr[29,8] = '(SP-RZK)'
$arr[29,9] = '3:07 PM'
$arr[29,11] = '0 hours, 2 minutes'
$arr[30,0] = 565
$arr[30,1] = 'Sunday, Jan 15'
$arr[30,2] = '3:35 PM'
$arr[30,3] = 'FR1543'
$arr[30,4] = 'Venice'
$arr[30,5] = '(TSF)'
$arr[30,6] = 'Buzz '
$arr[30,7] = 'B38M'
$arr[30,8] = '(SP-RZH)'
$arr[30,9] = '3:43 PM'
$arr[30,11] = '0 hours, 8 minutes'
$arr[31,0] = 566
$arr[31,1] = 'Sunday, Jan 15'
$arr[31,2] = '3:45 PM'
$arr[31,3] = 'FR3364'
$arr[31,4] = 'Berlin'
$arr[31,5] = '(BER)'
$arr[31,6] = 'Ryanair '
$arr[31,7] = 'B738'
$arr[31,8] = '(9H-QES)'
$arr[31,9] = '3:24 PM'
$arr[31,11] = '0 hours, -21 minutes'
$arr[32,0] = 567
$arr[32,1] = 'Sunday, Jan 15'
$arr[32,2] = '3:55 PM'
$arr[32,3] = 'LO3919'
$arr[32,4] = 'Warsaw'
$arr[32,5] = '(WAW)'
$arr[32,6] = 'LOT '
$arr[32,7] = 'E75S'
$arr[32,8] = '(SP-LIA)'
$arr[32,9] = '3:46 PM'
$arr[32,11] = '0 hours, -9 minutes'
$arr[33,0] = 568
$arr[33,1] = 'Sunday, Jan 15'
$arr[33,2] = '4:10 PM'
$arr[33,3] = 'KL1995'
$arr[33,4] = 'Amsterdam'
$arr[33,5] = '(AMS)'
$arr[33,6] = 'KLM '
$arr[33,7] = 'E190'
$arr[33,8] = '(PH-EXC)'
$arr[33,9] = '4:01 PM'
$arr[33,11] = '0 hours, -9 minutes'
$arr[34,0] = 569
$arr[34,1] = 'Sunday, Jan 15'
$arr[34,2] = '4:10 PM'
$arr[34,3] = 'U22113'
$arr[34,4] = 'London'
$arr[34,5] = '(LTN)'
$arr[34,6] = 'easyJet '
$arr[34,7] = 'A320'
$arr[34,8] = '(G-EZWX)'
$arr[34,9] = '3:41 PM'
$arr[34,11] = '0 hours, -29 minutes'
$arr[35,0] = 570
$arr[35,1] = 'Sunday, Jan 15'
$arr[35,2] = '4:15 PM'
$arr[35,3] = 'FR83'
$arr[35,4] = 'Eilat'
$arr[35,5] = '(ETM)'
$arr[35,6] = 'Buzz '
$arr[35,7] = 'B38M'
$arr[35,8] = '(SP-RZB)'
$arr[35,9] = '3:57 PM'
$arr[35,11] = '0 hours, -18 minutes'
$arr[36,0] = 571
$arr[36,1] = 'Sunday, Jan 15'
$arr[36,2] = '4:20 PM'
$arr[36,3] = 'W65074'
$arr[36,4] = 'Malaga'
$arr[36,5] = '(AGP)'
$arr[36,6] = 'Wizz Air '
$arr[36,7] = 'A21N'
$arr[36,8] = '(HA-LVO)'
$arr[36,9] = '3:50 PM'
$arr[36,11] = '0 hours, -30 minutes'

$ws.Range("A536:M572").Value = $arr
